$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 4 (old rows 4,5,6 shift down to 5,6,7).
$ws.Rows("4:4").Insert()

# Seed the new row 4 with a duplicate of row 3 (same PFR/accession/etc. block),
# then overwrite just the Sequence column with the new test sequence.
$ws.Range("A3:X3").Copy()
$ws.Range("A4").PasteSpecial()
$ws.Range("E4").Value = "XTESALSYAALILADSEIEISSEKLLTLTNAANVPVENIWADIFAKALDGQNLKDLLVNFSAGAAAPAGVAGGVAGGEAGEAEAEKEEEEAKEESDDDMGFGLFD"

# Match the saved selection/view state from the authored workbook.
$ws.Range("E4").Select() | Out-Null
